# Applies the data edits described by the commit "edited functions to be better"
# to Sheet1 of the active workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- New column header: I1 = "model" ----
$ws.Range("I1").Value = "model"

# ---- Row 8 (brackish water desalination): clear B8/C8 ----
$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()

# ---- Row 9 (seawater desalination): update values, add I9 ----
$ws.Range("B9").Value = 9.3423
$ws.Range("C9").Value = 0.7177
$ws.Range("E9").Value = 2.9129
$ws.Range("F9").Value = 0.6484
$ws.Range("I9").Value = 3

# ---- Row 11 (coagulation): add B11:G11 and I11 ----
$ws.Range("B11").Value = 0.222
$ws.Range("C11").Value = 1.516
$ws.Range("D11").Value = 3.071
$ws.Range("E11").Value = 0.347
$ws.Range("F11").Value = 1.448
$ws.Range("G11").Value = 2.726
$ws.Range("I11").Value = 1

# ---- New column values for existing rows ----
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("I5").Value = 2
$ws.Range("I6").Value = 2
$ws.Range("I7").Value = 2

# ---- Row 13 (new row): nanofiltration ----
$ws.Range("A13").Value = "nanofiltration"
$ws.Range("B13").Value = 7.14
$ws.Range("C13").Value = -0.22
$ws.Range("E13").Value = 0.44
$ws.Range("F13").Value = -0.13
$ws.Range("I13").Value = 2

# ---- Selection change recorded in the diff ----
$ws.Range("E16").Select()

$wb.Save()
